$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows above the current row 279 (shifts existing 279:297 down to 282:300,
# carrying formatting from the row above as Excel normally does on Rows.Insert).
$ws.Rows("279:281").Insert()

# Populate the 3 newly inserted rows with the new "Angeleno" weekly entries
# (Mercado/Región/Tipo/Producto/Categoría/Variedad columns match the rest of the table).

# Row 279: Angeleno / Especial
$ws.Cells.Item(279, 1).Value = 8
$ws.Cells.Item(279, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(279, 3).Value = "Coquimbo"
$ws.Cells.Item(279, 4).Value = 44714
$ws.Cells.Item(279, 5).Value = 4
$ws.Cells.Item(279, 6).Value = "Fruta"
$ws.Cells.Item(279, 7).Value = 100103
$ws.Cells.Item(279, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(279, 9).Value = 100103002
$ws.Cells.Item(279, 10).Value = "Ciruela"
$ws.Cells.Item(279, 11).Value = "Angeleno"
$ws.Cells.Item(279, 12).Value = "Especial"
$ws.Cells.Item(279, 13).Value = 16
$ws.Cells.Item(279, 14).Value = 230000
$ws.Cells.Item(279, 15).Value = 240000
$ws.Cells.Item(279, 16).Value = 235000
$ws.Cells.Item(279, 17).Value = "$/bins (450 kilos)"
$ws.Cells.Item(279, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(279, 19).Value = 522
$ws.Cells.Item(279, 20).Value = 450

# Row 280: Angeleno / Primera
$ws.Cells.Item(280, 1).Value = 8
$ws.Cells.Item(280, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(280, 3).Value = "Coquimbo"
$ws.Cells.Item(280, 4).Value = 44714
$ws.Cells.Item(280, 5).Value = 4
$ws.Cells.Item(280, 6).Value = "Fruta"
$ws.Cells.Item(280, 7).Value = 100103
$ws.Cells.Item(280, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(280, 9).Value = 100103002
$ws.Cells.Item(280, 10).Value = "Ciruela"
$ws.Cells.Item(280, 11).Value = "Angeleno"
$ws.Cells.Item(280, 12).Value = "Primera"
$ws.Cells.Item(280, 13).Value = 16
$ws.Cells.Item(280, 14).Value = 200000
$ws.Cells.Item(280, 15).Value = 210000
$ws.Cells.Item(280, 16).Value = 205000
$ws.Cells.Item(280, 17).Value = "$/bins (450 kilos)"
$ws.Cells.Item(280, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(280, 19).Value = 456
$ws.Cells.Item(280, 20).Value = 450

# Row 281: Angeleno / Segunda
$ws.Cells.Item(281, 1).Value = 8
$ws.Cells.Item(281, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(281, 3).Value = "Coquimbo"
$ws.Cells.Item(281, 4).Value = 44714
$ws.Cells.Item(281, 5).Value = 4
$ws.Cells.Item(281, 6).Value = "Fruta"
$ws.Cells.Item(281, 7).Value = 100103
$ws.Cells.Item(281, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(281, 9).Value = 100103002
$ws.Cells.Item(281, 10).Value = "Ciruela"
$ws.Cells.Item(281, 11).Value = "Angeleno"
$ws.Cells.Item(281, 12).Value = "Segunda"
$ws.Cells.Item(281, 13).Value = 16
$ws.Cells.Item(281, 14).Value = 170000
$ws.Cells.Item(281, 15).Value = 180000
$ws.Cells.Item(281, 16).Value = 175000
$ws.Cells.Item(281, 17).Value = "$/bins (450 kilos)"
$ws.Cells.Item(281, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(281, 19).Value = 389
$ws.Cells.Item(281, 20).Value = 450
